# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.784.72"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "2.499.04"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'588.31"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'176.23"
$ws.Range("E6").Value = "  +4.15%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D12").Value = "'4.94"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "'25.84"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "2.917.23"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000173"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.540.57"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "2.474.26"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D20").Value = "'352.65"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'70.65"
$ws.Range("E23").Value = "  +3.27%  "
$ws.Range("D24").Value = "'4.32"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").Value = "'1.76"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").Value = "2.624.53"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "'0.989"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "0.0₃0915"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").Value = "'510.00"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("D33").Value = "'1.79"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +8.59%  "
$ws.Range("D36").Value = "'163.91"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").Value = "'18.68"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("D44").Value = "'2.42"
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").Value = "'146.08"
$ws.Range("E45").Value = "  +3.67%  "
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").Value = "0.0₆0259"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "'0.587"
